$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value2 = "ECs"
$ws.Range("B2").Value2 = "Gnai2"
$ws.Range("C2").Value2 = "Adcy1"
$ws.Range("D2").Value2 = "ECs"
$ws.Range("E2").Value2 = 3
$ws.Range("F2").Value2 = 1
$ws.Range("G2").Value2 = 201.4397426666667
$ws.Range("H2").Value2 = 604.3192280000001
$ws.Range("I2").Value2 = 0.4833500233086392
$ws.Range("J2").Value2 = 0.4833500233086393
$ws.Range("K2").Value2 = 1
$ws.Range("L2").Value2 = 0.3333333333333333
$ws.Range("M2").Value2 = 0.04069533333333333
$ws.Range("N2").Value2 = 0.122086
$ws.Range("O2").Value2 = 0.1933524121143397
$ws.Range("P2").Value2 = 0.1933524121143397
$ws.Range("Q2").Value2 = 8.19765747440089
$ws.Range("R2").Value2 = 73.778917269608
$ws.Range("S2").Value2 = 0.0934568929022477
$ws.Range("T2").Value2 = 0.0934568929022477
$ws.Range("A3").Value2 = "ECs"
$ws.Range("B3").Value2 = "Gnai2"
$ws.Range("C3").Value2 = "Adcy1"
$ws.Range("D3").Value2 = "MuSCs"
$ws.Range("E3").Value2 = 3
$ws.Range("F3").Value2 = 1
$ws.Range("G3").Value2 = 201.4397426666667
$ws.Range("H3").Value2 = 604.3192280000001
$ws.Range("I3").Value2 = 0.4833500233086392
$ws.Range("J3").Value2 = 0.4833500233086393
$ws.Range("K3").Value2 = 2
$ws.Range("L3").Value2 = 0.6666666666666666
$ws.Range("M3").Value2 = 0.169777
$ws.Range("N3").Value2 = 0.509331
$ws.Range("O3").Value2 = 0.8066475878856604
$ws.Range("P3").Value2 = 0.8066475878856603
$ws.Range("Q3").Value2 = 34.19983519071867
$ws.Range("R3").Value2 = 307.798516716468
$ws.Range("S3").Value2 = 0.3898931304063916
$ws.Range("T3").Value2 = 0.3898931304063916
$ws.Range("A4").Value2 = "FAPs"
$ws.Range("B4").Value2 = "Gnai2"
$ws.Range("C4").Value2 = "Adcy1"
$ws.Range("D4").Value2 = "ECs"
$ws.Range("E4").Value2 = 3
$ws.Range("F4").Value2 = 1
$ws.Range("G4").Value2 = 65.41736466666667
$ws.Range("H4").Value2 = 196.252094
$ws.Range("I4").Value2 = 0.1569674599353791
$ws.Range("J4").Value2 = 0.1569674599353792
$ws.Range("K4").Value2 = 1
$ws.Range("L4").Value2 = 0.3333333333333333
$ws.Range("M4").Value2 = 0.04069533333333333
$ws.Range("N4").Value2 = 0.122086
$ws.Range("O4").Value2 = 0.1933524121143397
$ws.Range("P4").Value2 = 0.1933524121143397
$ws.Range("Q4").Value2 = 2.662181460898223
$ws.Range("R4").Value2 = 23.959633148084
$ws.Range("S4").Value2 = 0.03035003700196653
$ws.Range("T4").Value2 = 0.03035003700196653
$ws.Range("A5").Value2 = "FAPs"
$ws.Range("B5").Value2 = "Gnai2"
$ws.Range("C5").Value2 = "Adcy1"
$ws.Range("D5").Value2 = "MuSCs"
$ws.Range("E5").Value2 = 3
$ws.Range("F5").Value2 = 1
$ws.Range("G5").Value2 = 65.41736466666667
$ws.Range("H5").Value2 = 196.252094
$ws.Range("I5").Value2 = 0.1569674599353791
$ws.Range("J5").Value2 = 0.1569674599353792
$ws.Range("K5").Value2 = 2
$ws.Range("L5").Value2 = 0.6666666666666666
$ws.Range("M5").Value2 = 0.169777
$ws.Range("N5").Value2 = 0.509331
$ws.Range("O5").Value2 = 0.8066475878856604
$ws.Range("P5").Value2 = 0.8066475878856603
$ws.Range("Q5").Value2 = 11.10636392101267
$ws.Range("R5").Value2 = 99.957275289114
$ws.Range("S5").Value2 = 0.1266174229334126
$ws.Range("T5").Value2 = 0.1266174229334126
$ws.Range("A6").Value2 = "MuSCs"
$ws.Range("B6").Value2 = "Gnai2"
$ws.Range("C6").Value2 = "Adcy1"
$ws.Range("D6").Value2 = "ECs"
$ws.Range("E6").Value2 = 3
$ws.Range("F6").Value2 = 1
$ws.Range("G6").Value2 = 60.43484133333334
$ws.Range("H6").Value2 = 181.304524
$ws.Range("I6").Value2 = 0.1450120099461104
$ws.Range("J6").Value2 = 0.1450120099461104
$ws.Range("K6").Value2 = 1
$ws.Range("L6").Value2 = 0.3333333333333333
$ws.Range("M6").Value2 = 0.04069533333333333
$ws.Range("N6").Value2 = 0.122086
$ws.Range("O6").Value2 = 0.1933524121143397
$ws.Range("P6").Value2 = 0.1933524121143397
$ws.Range("Q6").Value2 = 2.459416013007111
$ws.Range("R6").Value2 = 22.134744117064
$ws.Range("S6").Value2 = 0.02803842190862905
$ws.Range("T6").Value2 = 0.02803842190862905
$ws.Range("A7").Value2 = "MuSCs"
$ws.Range("B7").Value2 = "Gnai2"
$ws.Range("C7").Value2 = "Adcy1"
$ws.Range("D7").Value2 = "MuSCs"
$ws.Range("E7").Value2 = 3
$ws.Range("F7").Value2 = 1
$ws.Range("G7").Value2 = 60.43484133333334
$ws.Range("H7").Value2 = 181.304524
$ws.Range("I7").Value2 = 0.1450120099461104
$ws.Range("J7").Value2 = 0.1450120099461104
$ws.Range("K7").Value2 = 2
$ws.Range("L7").Value2 = 0.6666666666666666
$ws.Range("M7").Value2 = 0.169777
$ws.Range("N7").Value2 = 0.509331
$ws.Range("O7").Value2 = 0.8066475878856604
$ws.Range("P7").Value2 = 0.8066475878856603
$ws.Range("Q7").Value2 = 10.26044605704933
$ws.Range("R7").Value2 = 92.34401451344401
$ws.Range("S7").Value2 = 0.1169735880374813
$ws.Range("T7").Value2 = 0.1169735880374813
$ws.Range("A8").Value2 = "Resolving-Mac"
$ws.Range("B8").Value2 = "Gnai2"
$ws.Range("C8").Value2 = "Adcy1"
$ws.Range("D8").Value2 = "ECs"
$ws.Range("E8").Value2 = 3
$ws.Range("F8").Value2 = 1
$ws.Range("G8").Value2 = 89.46554166666668
$ws.Range("H8").Value2 = 268.396625
$ws.Range("I8").Value2 = 0.2146705068098712
$ws.Range("J8").Value2 = 0.2146705068098712
$ws.Range("K8").Value2 = 1
$ws.Range("L8").Value2 = 0.3333333333333333
$ws.Range("M8").Value2 = 0.04069533333333333
$ws.Range("N8").Value2 = 0.122086
$ws.Range("O8").Value2 = 0.1933524121143397
$ws.Range("P8").Value2 = 0.1933524121143397
$ws.Range("Q8").Value2 = 3.640830039972223
$ws.Range("R8").Value2 = 32.76747035975001
$ws.Range("S8").Value2 = 0.04150706030149638
$ws.Range("T8").Value2 = 0.04150706030149638
$ws.Range("A9").Value2 = "Resolving-Mac"
$ws.Range("B9").Value2 = "Gnai2"
$ws.Range("C9").Value2 = "Adcy1"
$ws.Range("D9").Value2 = "MuSCs"
$ws.Range("E9").Value2 = 3
$ws.Range("F9").Value2 = 1
$ws.Range("G9").Value2 = 89.46554166666668
$ws.Range("H9").Value2 = 268.396625
$ws.Range("I9").Value2 = 0.2146705068098712
$ws.Range("J9").Value2 = 0.2146705068098712
$ws.Range("K9").Value2 = 2
$ws.Range("L9").Value2 = 0.6666666666666666
$ws.Range("M9").Value2 = 0.169777
$ws.Range("N9").Value2 = 0.509331
$ws.Range("O9").Value2 = 0.8066475878856604
$ws.Range("P9").Value2 = 0.8066475878856603
$ws.Range("Q9").Value2 = 15.18919126754167
$ws.Range("R9").Value2 = 136.702721407875
$ws.Range("S9").Value2 = 0.1731634465083748
$ws.Range("T9").Value2 = 0.1731634465083748
